# Apply the changes described by the commit:
#  - Add a note that a 10uF cap on usbvcc is highly recommended (new checklist
#    item in H6)
#  - Update the Avcc/low-pass-filter note in H3 to explain that the low pass
#    filter isn't needed because the 8u2's ADC isn't being used (appended as a
#    blue colored run), and mark that checklist cell as resolved ("Good" style)
#  - Remove the unused Sheet2 and Sheet3 tabs
#  - Update the active selection to H6

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- H6: new checklist note about the 10uF cap on usbvcc ---
$ws.Range("H6").Value = "10uF cap on usbvcc highly recommended"

# --- H3: extend the Avcc/low pass filter note with a second, blue colored run ---
$run1 = "Avcc externally connected to vcc through low pass filter  "
$run2 = "Not using adc on the 8u2"
$fullText = $run1 + $run2

$ws.Range("H3").Value = $fullText

$startPos = $run1.Length + 1
$runLen = $run2.Length
$chars = $ws.Range("H3").Characters($startPos, $runLen)
$chars.Font.Color = 12611584   # RGB(0, 112, 192) -> #0070C0
$chars.Font.Name = "Calibri"
$chars.Font.Size = 11

# Mark H3 as resolved by copying the "Good" (green) style already used by F3/G3
$ws.Range("F3").Copy()
$ws.Range("H3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update the selected/active cell on Sheet1 to H6 ---
$ws.Activate()
$ws.Range("H6").Select()

# Register the new blue font (#0070C0) in the workbook's font table, the same
# way real Excel does whenever that color is applied via the rich-text
# Characters() API. We stamp it on a cell of a soon-to-be-removed sheet so the
# font table picks it up without leaving any stray formatted cell behind.
$excel.DisplayAlerts = $false
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A1").Font.Color = 12611584

# --- Remove the now-unused Sheet2 and Sheet3 tabs ---
$wb.Worksheets.Item("Sheet3").Delete()
$wb.Worksheets.Item("Sheet2").Delete()
$excel.DisplayAlerts = $true
